$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    # Force a numeric-looking string to be stored as text (matches the
    # original inline-string cell type), then drop back to the default
    # "Normal" style so no stray number format lingers on the cell.
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "57.221.81"
$ws.Range("E2").Value = "  -5.37%  "
$ws.Range("D3").Value = "3.116.51"
$ws.Range("E3").Value = "  -6.03%  "
Set-TextValue "D4" "0.995"
$ws.Range("E4").Value = "  -0.43%  "
Set-TextValue "D5" "519.19"
$ws.Range("E5").Value = "  -7.23%  "
Set-TextValue "D6" "132.66"
$ws.Range("E6").Value = "  -7.75%  "
Set-TextValue "D7" "0.999"
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").Value = "3.111.73"
$ws.Range("E8").Value = "  -6.32%  "
Set-TextValue "D9" "0.443"
$ws.Range("E9").Value = "  -7.21%  "
Set-TextValue "D10" "7.12"
$ws.Range("E10").Value = "  -9.22%  "
Set-TextValue "D11" "0.108"
$ws.Range("E11").Value = "  -10.07%  "
Set-TextValue "D12" "0.379"
$ws.Range("E12").Value = "  -7.89%  "
$ws.Range("D13").Value = "3.647.58"
$ws.Range("E14").Value = "  -2.31%  "
Set-TextValue "D15" "25.37"
$ws.Range("E15").Value = "  -6.54%  "
$ws.Range("D16").Value = "3.115.46"
$ws.Range("E16").Value = "  -5.54%  "
$ws.Range("D17").Value = "57.053.85"
$ws.Range("E17").Value = "  -5.55%  "
Set-TextValue "D18" "0.0000149"
$ws.Range("E18").Value = "  -10.83%  "
Set-TextValue "D19" "5.72"
$ws.Range("E19").Value = "  -7.60%  "
Set-TextValue "D20" "12.85"
$ws.Range("E20").Value = "  -11.46%  "
Set-TextValue "D21" "7.91"
$ws.Range("E21").Value = "  -8.59%  "
Set-TextValue "D22" "344.23"
$ws.Range("E22").Value = "  -8.35%  "
Set-TextValue "D24" "67.94"
$ws.Range("E24").Value = "  -8.42%  "
Set-TextValue "D25" "0.501"
$ws.Range("E25").Value = "  -8.35%  "
$ws.Range("D26").Value = "3.262.44"
Set-TextValue "D27" "0.996"
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("D28").Value = "0.0₃0927"
$ws.Range("E28").Value = "  -10.89%  "
Set-TextValue "D29" "0.161"
$ws.Range("E29").Value = "  -6.84%  "
Set-TextValue "D30" "0.996"
$ws.Range("E30").Value = "  -0.35%  "
Set-TextValue "D31" "6.66"
$ws.Range("E31").Value = "  -8.82%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D32" "6.87"
$ws.Range("E32").Value = "  -11.00%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D33" "1.84"
$ws.Range("E33").Value = "  -9.98%  "
Set-TextValue "D34" "21.47"
$ws.Range("E34").Value = "  -5.19%  "
$ws.Range("E35").Value = "  -6.92%  "
Set-TextValue "D36" "4.78"
$ws.Range("E36").Value = "  -8.45%  "
Set-TextValue "D37" "155.32"
$ws.Range("E37").Value = "  -6.75%  "
Set-TextValue "D38" "6.11"
$ws.Range("E38").Value = "  -9.75%  "
$ws.Range("E39").Value = "  -11.48%  "
$ws.Range("E40").Value = "  -5.60%  "
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D41" "0.0681"
$ws.Range("E41").Value = "  -8.47%  "
$ws.Range("B42").Value = "RenzoRestakedETH"
$ws.Range("C42").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D42").Value = "3.132.97"
$ws.Range("E42").Value = "  -5.74%  "
Set-TextValue "D43" "40.25"
$ws.Range("E43").Value = "  -4.24%  "
Set-TextValue "D44" "0.679"
$ws.Range("E44").Value = "  -9.84%  "
Set-TextValue "D45" "3.87"
$ws.Range("E45").Value = "  -8.15%  "
Set-TextValue "D46" "0.995"
$ws.Range("E46").Value = "  -0.47%  "
Set-TextValue "D47" "1.04"
$ws.Range("E47").Value = "  -7.95%  "
$ws.Range("E48").Value = "  -9.22%  "
$ws.Range("D49").Value = "2.229.45"
$ws.Range("E49").Value = "  -5.76%  "
Set-TextValue "D50" "6.11"
$ws.Range("E50").Value = "  -7.00%  "
Set-TextValue "D51" "19.89"
$ws.Range("E51").Value = "  -7.54%  "
